$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the TabName in A2 from "CasesTab" to "ParticipantsTab"
$ws.Range("A2").Value = "ParticipantsTab"

# Select the edited cell, matching the resulting selection in the file
$ws.Range("A2").Select()
